$wb = $excel.ActiveWorkbook

# --- Sheet ED2A ---
$ws1 = $wb.Worksheets.Item("ED2A")
$ws1.Range("D1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)
$ws1.Range("E1").Value = "Plates"
$ws1.Range("E2").Value = "['BR00116991', 'BR00116992', 'BR00116993', 'BR00116994', 'BR00117015', 'BR00117016', 'BR00117017', 'BR00117019']"
$ws1.Range("E3").Value = "['BR00116991', 'BR00116992', 'BR00116993', 'BR00116994', 'BR00117020', 'BR00117021']"
$ws1.Range("E4").Value = "['BR00116991', 'BR00116992', 'BR00116993', 'BR00116994', 'BR00118050', 'BR00117006']"
$ws1.Range("E5").Value = "['BR00117015', 'BR00117016', 'BR00117017', 'BR00117019', 'BR00117020', 'BR00117021']"
$ws1.Range("E6").Value = "['BR00117015', 'BR00117016', 'BR00117017', 'BR00117019', 'BR00118050', 'BR00117006']"
$ws1.Range("E7").Value = "['BR00116991', 'BR00116992', 'BR00116993', 'BR00116994', 'BR00118041', 'BR00118042', 'BR00118043', 'BR00118044']"
$ws1.Range("E8").Value = "['BR00116991', 'BR00116992', 'BR00116993', 'BR00116994', 'BR00117000', 'BR00117003', 'BR00117004', 'BR00117005']"
$ws1.Range("E9").Value = "['BR00117015', 'BR00117016', 'BR00117017', 'BR00117019', 'BR00118041', 'BR00118042', 'BR00118043', 'BR00118044']"
$ws1.Range("E10").Value = "['BR00117015', 'BR00117016', 'BR00117017', 'BR00117019', 'BR00117000', 'BR00117003', 'BR00117004', 'BR00117005']"
$ws1.Range("E11").Value = "['BR00117020', 'BR00117021', 'BR00118050', 'BR00117006']"
$ws1.Range("E12").Value = "['BR00117020', 'BR00117021', 'BR00118041', 'BR00118042', 'BR00118043', 'BR00118044']"
$ws1.Range("E13").Value = "['BR00117020', 'BR00117021', 'BR00117000', 'BR00117003', 'BR00117004', 'BR00117005']"
$ws1.Range("E14").Value = "['BR00118050', 'BR00117006', 'BR00118041', 'BR00118042', 'BR00118043', 'BR00118044']"
$ws1.Range("E15").Value = "['BR00118050', 'BR00117006', 'BR00117000', 'BR00117003', 'BR00117004', 'BR00117005']"
$ws1.Range("E16").Value = "['BR00118041', 'BR00118042', 'BR00118043', 'BR00118044', 'BR00117000', 'BR00117003', 'BR00117004', 'BR00117005']"
$ws1.Range("E17").Value = "['BR00116991', 'BR00116992', 'BR00116993', 'BR00116994', 'BR00117008', 'BR00117009']"
$ws1.Range("E18").Value = "['BR00117015', 'BR00117016', 'BR00117017', 'BR00117019', 'BR00117008', 'BR00117009']"
$ws1.Range("E19").Value = "['BR00117020', 'BR00117021', 'BR00117008', 'BR00117009']"
$ws1.Range("E20").Value = "['BR00118050', 'BR00117006', 'BR00117008', 'BR00117009']"
$ws1.Range("E21").Value = "['BR00118041', 'BR00118042', 'BR00118043', 'BR00118044', 'BR00117008', 'BR00117009']"
$ws1.Range("E22").Value = "['BR00117000', 'BR00117003', 'BR00117004', 'BR00117005', 'BR00117008', 'BR00117009']"
$ws1.Range("E23").Value = "['BR00116991', 'BR00116992', 'BR00116993', 'BR00116994', 'BR00117054', 'BR00117055']"
$ws1.Range("E24").Value = "['BR00117015', 'BR00117016', 'BR00117017', 'BR00117019', 'BR00117054', 'BR00117055']"
$ws1.Range("E25").Value = "['BR00117020', 'BR00117021', 'BR00117054', 'BR00117055']"
$ws1.Range("E26").Value = "['BR00118050', 'BR00117006', 'BR00117054', 'BR00117055']"
$ws1.Range("E27").Value = "['BR00118041', 'BR00118042', 'BR00118043', 'BR00118044', 'BR00117054', 'BR00117055']"
$ws1.Range("E28").Value = "['BR00117000', 'BR00117003', 'BR00117004', 'BR00117005', 'BR00117054', 'BR00117055']"
$ws1.Range("E29").Value = "['BR00117008', 'BR00117009', 'BR00117054', 'BR00117055']"

# --- Sheet ED2B ---
$ws2 = $wb.Worksheets.Item("ED2B")
$ws2.Range("D1").Copy()
$ws2.Range("E1").PasteSpecial(-4122)
$ws2.Range("E1").Value = "Plates"
$ws2.Range("E2").Value = "['BR00116991', 'BR00116992', 'BR00116993', 'BR00116994', 'BR00117015', 'BR00117016', 'BR00117017', 'BR00117019']"
$ws2.Range("E3").Value = "['BR00116991', 'BR00116992', 'BR00116993', 'BR00116994', 'BR00117020', 'BR00117021']"
$ws2.Range("E4").Value = "['BR00116991', 'BR00116992', 'BR00116993', 'BR00116994', 'BR00118050', 'BR00117006']"
$ws2.Range("E5").Value = "['BR00117015', 'BR00117016', 'BR00117017', 'BR00117019', 'BR00117020', 'BR00117021']"
$ws2.Range("E6").Value = "['BR00117015', 'BR00117016', 'BR00117017', 'BR00117019', 'BR00118050', 'BR00117006']"
$ws2.Range("E7").Value = "['BR00117020', 'BR00117021', 'BR00118050', 'BR00117006']"
$ws2.Range("E8").Value = "['BR00116991', 'BR00116992', 'BR00116993', 'BR00116994', 'BR00117050', 'BR00117051', 'BR00117052', 'BR00117053']"
$ws2.Range("E9").Value = "['BR00117015', 'BR00117016', 'BR00117017', 'BR00117019', 'BR00117050', 'BR00117051', 'BR00117052', 'BR00117053']"
$ws2.Range("E10").Value = "['BR00116991', 'BR00116992', 'BR00116993', 'BR00116994', 'BR00118041', 'BR00118042', 'BR00118043', 'BR00118044']"
$ws2.Range("E11").Value = "['BR00116991', 'BR00116992', 'BR00116993', 'BR00116994', 'BR00117000', 'BR00117003', 'BR00117004', 'BR00117005']"
$ws2.Range("E12").Value = "['BR00117015', 'BR00117016', 'BR00117017', 'BR00117019', 'BR00118041', 'BR00118042', 'BR00118043', 'BR00118044']"
$ws2.Range("E13").Value = "['BR00117015', 'BR00117016', 'BR00117017', 'BR00117019', 'BR00117000', 'BR00117003', 'BR00117004', 'BR00117005']"
$ws2.Range("E14").Value = "['BR00117020', 'BR00117021', 'BR00117050', 'BR00117051', 'BR00117052', 'BR00117053']"
$ws2.Range("E15").Value = "['BR00118050', 'BR00117006', 'BR00117050', 'BR00117051', 'BR00117052', 'BR00117053']"
$ws2.Range("E16").Value = "['BR00117020', 'BR00117021', 'BR00118041', 'BR00118042', 'BR00118043', 'BR00118044']"
$ws2.Range("E17").Value = "['BR00117020', 'BR00117021', 'BR00117000', 'BR00117003', 'BR00117004', 'BR00117005']"
$ws2.Range("E18").Value = "['BR00118050', 'BR00117006', 'BR00118041', 'BR00118042', 'BR00118043', 'BR00118044']"
$ws2.Range("E19").Value = "['BR00118050', 'BR00117006', 'BR00117000', 'BR00117003', 'BR00117004', 'BR00117005']"
$ws2.Range("E20").Value = "['BR00117050', 'BR00117051', 'BR00117052', 'BR00117053', 'BR00118041', 'BR00118042', 'BR00118043', 'BR00118044']"
$ws2.Range("E21").Value = "['BR00117050', 'BR00117051', 'BR00117052', 'BR00117053', 'BR00117000', 'BR00117003', 'BR00117004', 'BR00117005']"
$ws2.Range("E22").Value = "['BR00118041', 'BR00118042', 'BR00118043', 'BR00118044', 'BR00117000', 'BR00117003', 'BR00117004', 'BR00117005']"
